$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4529.6
$ws.Range("I12").Value = 216
$ws.Range("J12").Value = 11000
$ws.Range("K12").Value = 216
$ws.Range("L12").Value = 11000
$ws.Range("M12").Value = -46
$ws.Range("N12").Value = -11340
$ws.Range("H33").Value = 264.2857
$ws.Range("H53").Value = 162.375
$ws.Range("I53").Value = 145.36363
$ws.Range("J53").Value = 199.8
$ws.Range("K53").Value = 145.36363
$ws.Range("L53").Value = 199.8
$ws.Range("M53").Value = 491.63637
$ws.Range("N53").Value = -1473.8
$ws.Range("H70").Value = 3445.1
$ws.Range("I70").Value = 2050.3333
$ws.Range("K70").Value = 6150.999899999999
$ws.Range("M70").Value = -5880.999899999999
$ws.Range("H73").Value = 3445.1
$ws.Range("I73").Value = 2050.3333
$ws.Range("K73").Value = 6150.999899999999
$ws.Range("M73").Value = -5214.999899999999
$ws.Range("H96").Value = 1892.0476
$ws.Range("I96").Value = 874.4
$ws.Range("K96").Value = 2623.2
$ws.Range("M96").Value = -1250.2
$ws.Range("H130").Value = 94640
$ws.Range("J130").Value = 94640
$ws.Range("L130").Value = 94640
$ws.Range("N130").Value = -104680
$ws.Range("H137").Value = 5707.2085
$ws.Range("I137").Value = 1612.409
$ws.Range("J137").Value = 50750
$ws.Range("K137").Value = 4837.227000000001
$ws.Range("L137").Value = 152250
$ws.Range("M137").Value = -2287.227000000001
$ws.Range("N137").Value = -157350

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 554
$ws.Range("I19").Value = 1008
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 1008
$ws.Range("L19").Value = 100
$ws.Range("M19").Value = -779
$ws.Range("N19").Value = -558
$ws.Range("H61").Value = 11701.091
$ws.Range("I61").Value = 14633.625
$ws.Range("J61").Value = 3881
$ws.Range("K61").Value = 14633.625
$ws.Range("L61").Value = 3881
$ws.Range("M61").Value = -14421.625
$ws.Range("N61").Value = -4305
$ws.Range("H102").Value = 1053.2222
$ws.Range("I102").Value = 903.625
$ws.Range("K102").Value = 903.625
$ws.Range("M102").Value = 718.375
$ws.Range("H122").Value = 1719.8572
$ws.Range("I122").Value = 1719.8572
$ws.Range("K122").Value = 5159.571599999999
$ws.Range("M122").Value = -2709.571599999999
$ws.Range("H136").Value = 11701.091
$ws.Range("I136").Value = 14633.625
$ws.Range("J136").Value = 3881
$ws.Range("K136").Value = 43900.875
$ws.Range("L136").Value = 11643
$ws.Range("M136").Value = -41350.875
$ws.Range("N136").Value = -16743

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5400.8
$ws.Range("I134").Value = 7054.857
$ws.Range("J134").Value = 1541.3334
$ws.Range("K134").Value = 21164.571
$ws.Range("L134").Value = 4624.0002
$ws.Range("M134").Value = -18629.571
$ws.Range("N134").Value = -9694.0002

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2493.5
$ws.Range("I22").Value = 1014.5833
$ws.Range("J22").Value = 4711.875
$ws.Range("K22").Value = 1014.5833
$ws.Range("L22").Value = 4711.875
$ws.Range("M22").Value = -664.5833
$ws.Range("N22").Value = -5411.875
$ws.Range("H25").Value = 6606.3335
$ws.Range("I25").Value = 6606.3335
$ws.Range("K25").Value = 6606.3335
$ws.Range("M25").Value = -6432.3335
$ws.Range("H86").Value = 42928.73
$ws.Range("I86").Value = 105595.75
$ws.Range("K86").Value = 105595.75
$ws.Range("M86").Value = -104472.75
$ws.Range("H89").Value = 42928.73
$ws.Range("I89").Value = 105595.75
$ws.Range("K89").Value = 527978.75
$ws.Range("M89").Value = -522362.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5270241.5
$ws.Range("I4").Value = 12454779
$ws.Range("K4").Value = 37364337
$ws.Range("M4").Value = -37364225
$ws.Range("H37").Value = 45453.637
$ws.Range("J37").Value = 45453.637
$ws.Range("L37").Value = 136360.911
$ws.Range("N37").Value = -136584.911
$ws.Range("H141").Value = 13842.857
$ws.Range("I141").Value = 5950
$ws.Range("J141").Value = 17000
$ws.Range("K141").Value = 17850
$ws.Range("L141").Value = 51000
$ws.Range("M141").Value = -12670
$ws.Range("N141").Value = -61360

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 111.882355
$ws.Range("I2").Value = 57.875
$ws.Range("J2").Value = 159.88889
$ws.Range("K2").Value = 57.875
$ws.Range("L2").Value = 159.88889
$ws.Range("M2").Value = 55.125
$ws.Range("N2").Value = -385.88889
$ws.Range("H21").Value = 14950
$ws.Range("I21").Value = 4800
$ws.Range("J21").Value = 18333.334
$ws.Range("K21").Value = 4800
$ws.Range("L21").Value = 18333.334
$ws.Range("M21").Value = -4627
$ws.Range("N21").Value = -18679.334
$ws.Range("H30").Value = 14950
$ws.Range("I30").Value = 4800
$ws.Range("J30").Value = 18333.334
$ws.Range("K30").Value = 4800
$ws.Range("L30").Value = 18333.334
$ws.Range("M30").Value = -4695
$ws.Range("N30").Value = -18543.334
$ws.Range("H80").Value = 1770
$ws.Range("J80").Value = 1980
$ws.Range("L80").Value = 1980
$ws.Range("N80").Value = -3976
$ws.Range("H83").Value = 1770
$ws.Range("J83").Value = 1980
$ws.Range("L83").Value = 9900
$ws.Range("N83").Value = -19884
$ws.Range("H102").Value = 4257.5
$ws.Range("I102").Value = 912.2
$ws.Range("K102").Value = 912.2
$ws.Range("M102").Value = 709.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 23374.75
$ws.Range("J59").Value = 23374.75
$ws.Range("L59").Value = 23374.75
$ws.Range("N59").Value = -24682.75
$ws.Range("H82").Value = 14999
$ws.Range("I82").Value = 19166.584
$ws.Range("J82").Value = 2496.25
$ws.Range("K82").Value = 19166.584
$ws.Range("L82").Value = 2496.25
$ws.Range("M82").Value = -18805.584
$ws.Range("N82").Value = -3218.25
$ws.Range("H85").Value = 14999
$ws.Range("I85").Value = 19166.584
$ws.Range("J85").Value = 2496.25
$ws.Range("K85").Value = 19166.584
$ws.Range("L85").Value = 2496.25
$ws.Range("M85").Value = -17918.584
$ws.Range("N85").Value = -4992.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2999
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2999
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2999
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -5745
$ws.Range("H100").Value = 1514.5
$ws.Range("I100").Value = 1096
$ws.Range("J100").Value = 1933
$ws.Range("K100").Value = 2192
$ws.Range("L100").Value = 3866
$ws.Range("M100").Value = -1651
$ws.Range("N100").Value = -4948
$ws.Range("H113").Value = 1568.262
$ws.Range("I113").Value = 1195.5
$ws.Range("K113").Value = 3586.5
$ws.Range("M113").Value = -1416.5
$ws.Range("H136").Value = 155605.14
$ws.Range("I136").Value = 14872.667
$ws.Range("K136").Value = 44618.001
$ws.Range("M136").Value = -42068.001
